$wb = $excel.ActiveWorkbook

# --- Sheet: "Trends Status" ---
$ws1 = $wb.Worksheets.Item("Trends Status")

$ws1.Range("C2").Value = 3
$ws1.Range("E2").Value = 20

$ws1.Range("D3").Value = 12.5
$ws1.Range("E3").Value = 6.7

$ws1.Range("B4").Value = 3
$ws1.Range("C4").Value = 6
$ws1.Range("D4").Value = 37.5
$ws1.Range("E4").Value = 40

$ws1.Range("C5").Value = 3
$ws1.Range("D5").Value = 12.5
$ws1.Range("E5").Value = 20

$ws1.Range("B6").Value = 3
$ws1.Range("D6").Value = 37.5
$ws1.Range("E6").Value = 13.3

$ws1.Range("B7").Value = 59
$ws1.Range("C7").Value = 65

# --- Sheet: "Species qualification" ---
$ws4 = $wb.Worksheets.Item("Species qualification")

$ws4.Range("C3").Value = 8
$ws4.Range("C4").Value = 15
